{"js": "// Replace the date line and each \"a\u00f7b=\" drill cell with the new values,\n// in document order (the same text can repeat, so we walk paragraphs in\n// order rather than doing a single global text replace).\nconst replacements = [\n  \"2024-07-20 Saturday\", \"2024-07-21 Sunday\",\n  \"13\u00f72=\", \"39\u00f75=\",\n  \"69\u00f72=\", \"19\u00f73=\",\n  \"60\u00f78=\", \"50\u00f78=\",\n  \"93\u00f76=\", \"52\u00f78=\",\n  \"38\u00f75=\", \"70\u00f76=\",\n  \"13\u00f73=\", \"18\u00f77=\",\n  \"97\u00f78=\", \"54\u00f76=\",\n  \"10\u00f78=\", \"29\u00f77=\",\n  \"92\u00f75=\", \"41\u00f72=\",\n  \"92\u00f75=\", \"74\u00f73=\",\n  \"99\u00f72=\", \"30\u00f79=\",\n  \"59\u00f73=\", \"53\u00f75=\",\n  \"45\u00f76=\", \"22\u00f75=\",\n  \"93\u00f72=\", \"61\u00f76=\",\n  \"34\u00f76=\", \"41\u00f73=\",\n  \"98\u00f79=\", \"71\u00f74=\",\n  \"82\u00f75=\", \"17\u00f73=\",\n  \"35\u00f79=\", \"48\u00f78=\",\n  \"16\u00f76=\", \"59\u00f72=\",\n  \"18\u00f78=\", \"27\u00f72=\",\n  \"80\u00f77=\", \"21\u00f78=\",\n  \"18\u00f73=\", \"93\u00f79=\",\n  \"63\u00f75=\", \"52\u00f72=\",\n  \"64\u00f76=\", \"44\u00f75=\",\n  \"54\u00f75=\", \"30\u00f72=\",\n];\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nlet ptr = 0;\nfor (const paragraph of paragraphs.items) {\n  if (ptr >= replacements.length) break;\n  const expected = replacements[ptr];\n  if (paragraph.text === expected) {\n    paragraph.getRange().insertText(replacements[ptr + 1], \"Replace\");\n    ptr += 2;\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the worksheet date and refresh every \"a\u00f7b=\" drill cell with new\n# operands, in document/table order.\n$d = $word.ActiveDocument\n\n$d.Content.Find.Execute(\n  \"2024-07-20 Saturday\", $false, $false, $false, $false, $false, $true,\n  1, $false, \"2024-07-21 Sunday\", 2\n)\n\n$newValues = @(\n  \"39\u00f75=\", \"19\u00f73=\", \"50\u00f78=\", \"52\u00f78=\", \"70\u00f76=\",\n  \"18\u00f77=\", \"54\u00f76=\", \"29\u00f77=\", \"41\u00f72=\", \"74\u00f73=\",\n  \"30\u00f79=\", \"53\u00f75=\", \"22\u00f75=\", \"61\u00f76=\", \"41\u00f73=\",\n  \"71\u00f74=\", \"17\u00f73=\", \"48\u00f78=\", \"59\u00f72=\", \"27\u00f72=\",\n  \"21\u00f78=\", \"93\u00f79=\", \"52\u00f72=\", \"44\u00f75=\", \"30\u00f72=\"\n)\n\n$dataRows = @(1, 5, 9, 13, 17)\n$t = $d.Tables.Item(1)\n\n$i = 0\nforeach ($r in $dataRows) {\n  for ($c = 1; $c -le 5; $c++) {\n    $cell = $t.Cell($r, $c)\n    $cell.Range.Text = $newValues[$i]\n    $i = $i + 1\n  }\n}\n"}
